$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Update the TCode value used for filtering (E2) from co03 to co02,
# since the ID column for filtering wouldn't always display properly.
$ws.Range("E2").Value = "co02"

# Update the saved selection on the Global sheet from B2 to A2,
# using a different navigation method to reach the ID filter, while
# preserving whichever sheet/tab was active before the edit.
$originalActiveSheet = $wb.ActiveSheet
$ws.Range("A2").Select()
$originalActiveSheet.Activate()
